$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet ARM (sheet2.xml)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 63
$ws.Range("H63").Value = 3104.6428
$ws.Range("I63").Value = 3158.8462
$ws.Range("J63").Value = 2400
$ws.Range("K63").Value = 3158.8462
$ws.Range("L63").Value = 2400
$ws.Range("M63").Value = -2472.8462
$ws.Range("N63").Value = -3772

# Row 66
$ws.Range("H66").Value = 3104.6428
$ws.Range("I66").Value = 3158.8462
$ws.Range("J66").Value = 2400
$ws.Range("K66").Value = 15794.231
$ws.Range("L66").Value = 12000
$ws.Range("M66").Value = -12362.231
$ws.Range("N66").Value = -18864

# ---------------------------------------------------------------------
# Sheet BSM (sheet3.xml)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 35
$ws.Range("H35").Value = 29559.8
$ws.Range("J35").Value = 29559.8
$ws.Range("L35").Value = 29559.8
$ws.Range("N35").Value = -30179.8

# Row 82
$ws.Range("H82").Value = 14337.667
$ws.Range("I82").Value = 4756.5
$ws.Range("J82").Value = 33500
$ws.Range("K82").Value = 4756.5
$ws.Range("L82").Value = 33500
$ws.Range("M82").Value = -4373.5
$ws.Range("N82").Value = -34266

# Row 85
$ws.Range("H85").Value = 14337.667
$ws.Range("I85").Value = 4756.5
$ws.Range("J85").Value = 33500
$ws.Range("K85").Value = 4756.5
$ws.Range("L85").Value = 33500
$ws.Range("M85").Value = -3430.5
$ws.Range("N85").Value = -36152

# Row 107
$ws.Range("H107").Value = 1545.4
$ws.Range("I107").Value = 1422.1428
$ws.Range("J107").Value = 1833
$ws.Range("K107").Value = 1422.1428
$ws.Range("L107").Value = 1833
$ws.Range("M107").Value = 497.8571999999999
$ws.Range("N107").Value = -5673

# Rows whose H:N values are cleared entirely (row 121 and 136 are left untouched)
$clearRows = @(117,118,119,120,122,123,124,125,126,127,128,129,130,131,132,133,134,135,137,138,139,140,141)
foreach ($r in $clearRows) {
    $ws.Range("H$r`:N$r").ClearContents()
}

# ---------------------------------------------------------------------
# Sheet GSM (sheet6.xml)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 117
$ws.Range("H117").Value = 30000
$ws.Range("J117").Value = 30000
$ws.Range("L117").Value = 30000
$ws.Range("N117").Value = -36884

# ---------------------------------------------------------------------
# Sheet LTW (sheet7.xml)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Rows whose H:N values are cleared entirely (row 126 is left untouched)
$clearRows2 = @(124,125,127,128,129,130,131,132,133,134,135,136,137,138,139,140,141)
foreach ($r in $clearRows2) {
    $ws.Range("H$r`:N$r").ClearContents()
}
